$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 3; $r++) {
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = "yolima"
    $ws.Cells.Item($r, 3).Value = "hola mundo"
    $ws.Cells.Item($r, 4).Value = "hola"
}

$ws.Range("A3:D3").Select() | Out-Null
